# Aufwandschätzung.xlsx - extend interval to support newer coin counter
# hardware: split the single "Stunden" effort column into a "Stunden min"
# and a "Stunden max" column with new effort estimations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# New minimum-hours values per row (2..31) - the former single "Stunden"
# column becomes the new minimum estimate.
$minValues = @{
    2  = 0.5
    3  = 0.5
    4  = 0.5
    5  = 1
    6  = 0.5
    7  = 0.5
    8  = 2
    9  = 1
    10 = 0.5
    11 = 1
    12 = 0.5
    13 = 0.5
    14 = 1
    15 = 0.5
    16 = 0.5
    17 = 1
    18 = 1
    19 = 0.5
    20 = 0.5
    21 = 0.5
    22 = 0.5
    23 = 0.5
    24 = 0.5
    25 = 0.5
    26 = 0.5
    27 = 0.5
    28 = 1
    29 = 1
    30 = 1
    31 = 1
}

# Move the old "Stunden" values (column C) over into the new column D
# ("Stunden max"), then overwrite column C with the new minimum estimate.
for ($r = 2; $r -le 31; $r++) {
    $oldHours = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value = $oldHours
    $ws.Cells.Item($r, 3).Value = $minValues[$r]
}

# Fix typo in the "Holz" row about the back panel.
$ws.Range("B25").Value = "Bohren, Dübeln, Leimen / Einpassen Rückwand inkl. Kabeldurchführung"

# Rename header columns: "Stunden" -> "Stunden" + newline + "min", and add
# a new "Stunden" + newline + "max" header in column D (outside the table).
$ws.Range("C1").Value = "Stunden" + $nl + "min"
$ws.Range("D1").Value = "Stunden" + $nl + "max"

# Match formatting of the other header cells (bold font) and enable
# word-wrap for the two-line headers.
$ws.Range("C1:D1").Font.Bold = $true
$ws.Range("C1:D1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Totals row: add SUM for the new max column and an average formula that
# also references the (table-qualified) min column.
$ws.Range("D33").Formula = "=SUM(D1:D32)"
$ws.Range("D33").Font.Bold = $true
$ws.Range("E33").Formula = "=(D33+Tabelle2[[#This Row],[Stunden" + $nl + "min]])/2"

# The multi-line formula text above makes the engine mis-estimate row 33's
# height; auto-fit it back down to the regular default height.
$ws.Rows.Item(33).AutoFit()

# Leave the selection on A33, matching the saved workbook state.
$ws.Range("A33").Select()
